$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a cell's value as literal text, preventing Excel's automatic
# number/date inference (e.g. "1.00" -> 1, "0.489" -> 0.489 numeric) while
# keeping the cell free of any extra style index (matches original unstyled
# data cells).
function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# Auto-generated cell updates
$updates = @(
    @{ Cell = "D2"; Value = "41.944.00" }
    @{ Cell = "E2"; Value = "  -0.21%  " }
    @{ Cell = "D3"; Value = "2.276.69" }
    @{ Cell = "E3"; Value = "  +0.89%  " }
    @{ Cell = "D4"; Value = "1.00" }
    @{ Cell = "E4"; Value = "  -0.04%  " }
    @{ Cell = "D5"; Value = "305.69" }
    @{ Cell = "E5"; Value = "  +1.28%  " }
    @{ Cell = "D6"; Value = "93.19" }
    @{ Cell = "E6"; Value = "  +0.43%  " }
    @{ Cell = "E7"; Value = "  -0.77%  " }
    @{ Cell = "D8"; Value = "1.00" }
    @{ Cell = "D9"; Value = "0.489" }
    @{ Cell = "E9"; Value = "  +1.16%  " }
    @{ Cell = "D10"; Value = "32.80" }
    @{ Cell = "E10"; Value = "  +0.26%  " }
    @{ Cell = "E11"; Value = "  -0.28%  " }
    @{ Cell = "E12"; Value = "  -1.82%  " }
    @{ Cell = "D13"; Value = "6.71" }
    @{ Cell = "E13"; Value = "  +0.43%  " }
    @{ Cell = "D14"; Value = "2.625.35" }
    @{ Cell = "E14"; Value = "  +0.69%  " }
    @{ Cell = "E15"; Value = "  +1.74%  " }
    @{ Cell = "D16"; Value = "2.275.33" }
    @{ Cell = "E16"; Value = "  +0.86%  " }
    @{ Cell = "D17"; Value = "0.784" }
    @{ Cell = "E17"; Value = "  +3.57%  " }
    @{ Cell = "D18"; Value = "41.870.66" }
    @{ Cell = "E18"; Value = "  -0.13%  " }
    @{ Cell = "D19"; Value = "12.87" }
    @{ Cell = "E19"; Value = "  +5.76%  " }
    @{ Cell = "E20"; Value = "  +1.09%  " }
    @{ Cell = "E21"; Value = "  +0.50%  " }
    @{ Cell = "D22"; Value = "68.01" }
    @{ Cell = "E22"; Value = "  +1.11%  " }
    @{ Cell = "D23"; Value = "244.06" }
    @{ Cell = "E23"; Value = "  +0.91%  " }
    @{ Cell = "D24"; Value = "2.61" }
    @{ Cell = "E24"; Value = "  +1.17%  " }
    @{ Cell = "E25"; Value = "  +1.92%  " }
    @{ Cell = "D27"; Value = "24.07" }
    @{ Cell = "E27"; Value = "  +0.38%  " }
    @{ Cell = "D28"; Value = "9.71" }
    @{ Cell = "E28"; Value = "  +0.28%  " }
    @{ Cell = "E29"; Value = "  -0.44%  " }
    @{ Cell = "D30"; Value = "34.98" }
    @{ Cell = "E30"; Value = "  +2.44%  " }
    @{ Cell = "D31"; Value = "159.15" }
    @{ Cell = "E31"; Value = "  +0.35%  " }
    @{ Cell = "D32"; Value = "5.39" }
    @{ Cell = "E32"; Value = "  +4.49%  " }
    @{ Cell = "E33"; Value = "  -0.01%  " }
    @{ Cell = "D34"; Value = "0.0744" }
    @{ Cell = "E34"; Value = "  +0.09%  " }
    @{ Cell = "D35"; Value = "3.04" }
    @{ Cell = "E35"; Value = "  -0.70%  " }
    @{ Cell = "D36"; Value = "17.39" }
    @{ Cell = "E36"; Value = "  +4.61%  " }
    @{ Cell = "E37"; Value = "  -1.29%  " }
    @{ Cell = "E38"; Value = "  +0.31%  " }
    @{ Cell = "E39"; Value = "  +0.54%  " }
    @{ Cell = "E40"; Value = "  -0.32%  " }
    @{ Cell = "D41"; Value = "3.95" }
    @{ Cell = "E41"; Value = "  +0.25%  " }
    @{ Cell = "B42"; Value = "EnergySwap" }
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" }
    @{ Cell = "D42"; Value = "19.80" }
    @{ Cell = "E42"; Value = "  -0.85%  " }
    @{ Cell = "B43"; Value = "Maker" }
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr" }
    @{ Cell = "D43"; Value = "2.009.66" }
    @{ Cell = "E43"; Value = "  -2.00%  " }
    @{ Cell = "D44"; Value = "2.26" }
    @{ Cell = "E44"; Value = "  +11.24%  " }
    @{ Cell = "E45"; Value = "  +1.18%  " }
    @{ Cell = "D46"; Value = "10.27" }
    @{ Cell = "E46"; Value = "  +1.54%  " }
    @{ Cell = "D47"; Value = "2.92" }
    @{ Cell = "E47"; Value = "  +0.61%  " }
    @{ Cell = "D48"; Value = "53.60" }
    @{ Cell = "E48"; Value = "  +3.25%  " }
    @{ Cell = "E49"; Value = "  +2.96%  " }
    @{ Cell = "B50"; Value = "TrustWalletToken" }
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt" }
    @{ Cell = "D50"; Value = "1.15" }
    @{ Cell = "E50"; Value = "  +0.60%  " }
    @{ Cell = "B51"; Value = "Stacks" }
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx" }
    @{ Cell = "D51"; Value = "1.51" }
    @{ Cell = "E51"; Value = "  -0.30%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    Set-TextValue $cell $u.Value
}
